$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New trade row appended to the bottom of the sheet
$ws.Range("A4").Value = 10031.92
$ws.Range("B4").Value = 9983
$ws.Range("C4").Value = 307.20999999999998
$ws.Range("D4").Value = 308.70999999999998
$ws.Range("E4").Value = $false
$ws.Range("F4").Value = 0.49
$ws.Range("G4").Value = 42609.503923611112
$ws.Range("G4").NumberFormat = "m/d/yy h:mm"
$ws.Range("H4").Value = $true

# Column A's "best fit" width grows by a hair now that the new number is in the mix
$ws.Columns.Item(1).ColumnWidth = 8.14
